$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Price updates (transport calculator surcharge increases) ---

# Rotterdam (POD) rows with MOTORCYCLE: price drops from 700 to 600
$rotterdamMotoRows = @(6, 11, 16, 21, 26, 31, 36)
foreach ($r in $rotterdamMotoRows) {
    $ws.Range("D$r").Value = 600
}

# Varna (POD) rows 37-71: prices change from a flat 700 to a
# per-vehicle-type tier (CAR/SUV/LARGE SUV/PICKUP = 400, MOTORCYCLE = 300),
# repeated for each of the 7 origin ports (POL).
$varnaStart = 37
for ($block = 0; $block -lt 7; $block++) {
    $base = $varnaStart + ($block * 5)
    $ws.Range("D$base").Value = 400
    $ws.Range("D$($base + 1)").Value = 400
    $ws.Range("D$($base + 2)").Value = 400
    $ws.Range("D$($base + 3)").Value = 400
    $ws.Range("D$($base + 4)").Value = 300
}

# --- View/selection state ---
# Scroll the sheet so row 31 is near the top, and leave the selection on H67
# (matches the author's last-saved cursor position).
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("H67").Select()
